# Data update for Germany
# Only Germany (DE) remains selected (1) in the Region_selection sheet;
# all other previously-selected European countries (rows 12-41, excluding
# the already-zero ES/HU rows and the DE row itself) are deselected (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Region_selection")

$rowsToClear = @(12,13,14,15,16,18,19,21,22,23,24,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41)

foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Update the view/selection on the sheet to match the saved state.
$ws.Activate()
$ws.Range("D27").Select()
